# Applies the commit "Se completan las métricas con reporte de code analyzer"
# - Fills in the remaining "Líneas Reales" (M column) values for increments 2-5
# - Adjusts the estimated/real execution time for increment 5 (row 21)
# - Adds a brand-new increment 6 (row 22): "EnvasadoraTest"
# - Fills in the previously-empty "Ejecución de la Prueba" timing row (row 30)
# - Moves the active selection to M23

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Métricas")

# --- Desarrollo y correctivos: Líneas Reales (M) for existing increments ---
$ws.Range("M18").Value = 11
$ws.Range("M19").Value = 38
$ws.Range("M20").Value = 20

# --- Increment 5 (row 21): updated estimated time & real end time ---
# (plain decimals below -- the host PowerShell parser chokes on "E-2"
# scientific-notation literals, so express the day-fractions as exact
# minute/second ratios instead)
$ws.Range("G21").Value = 40/1440
$ws.Range("I21").Value = (16*3600+25*60)/86400
$ws.Range("M21").Value = 83

# --- New increment 6 (row 22): "EnvasadoraTest" ---
$ws.Range("C22").Value = "EnvasadoraTest"
$ws.Range("F22").Value = 60
$ws.Range("G22").Value = 5/1440
$ws.Range("H22").Value = (16*3600+25*60)/86400
$ws.Range("I22").Value = 16.5/24
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 53

# --- Ejecución de la Prueba timing (row 30) ---
$ws.Range("B30").Value = 20/1440
$ws.Range("C30").Value = 16.5/24
$ws.Range("D30").Value = (16*3600+46*60)/86400

# --- Recalculate so every dependent formula (totals, chart cache, % of
#     total, etc.) picks up the new inputs ---
$excel.Calculate()

# --- Restore the active cell/selection shown when the file was saved ---
$ws.Range("M23").Select()
